# ReturnRental table refresh: the scraper re-ran and produced a new batch of
# rows. Row 8 (845218) is gone entirely, and the remaining rows 2-7 get new
# "Store", "Klant nummer" and "Serienummer" values (some serials are now
# alphanumeric strings instead of the long numeric IMEI-style values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Id 845212) ---
$ws.Range("H2").Value = "Sint-Denijs-Westrem"
$ws.Range("P2").Value = 1011782546
$ws.Range("V2").Value = 226415171533

# --- Row 3 (Id 845213) ---
$ws.Range("H3").Value = "Sint-Denijs-Westrem"
$ws.Range("P3").Value = 1011782546
$ws.Range("V3").Value = "DC537CF69DF4"

# --- Row 4 (Id 845214) ---
$ws.Range("H4").Value = "Sint-Denijs-Westrem"
$ws.Range("P4").Value = 693300224
$ws.Range("V4").Value = 965346592830

# --- Row 5 (Id 845215) ---
$ws.Range("H5").Value = "Sint-Denijs-Westrem"
$ws.Range("P5").Value = 693300224
$ws.Range("V5").Value = "AC220594A0BB"

# --- Row 6 (Id 845216) ---
$ws.Range("H6").Value = "Sint-Denijs-Westrem"
$ws.Range("P6").Value = 693300224
$ws.Range("V6").Value = "1907318070029335"

# --- Row 7 (Id 845217) ---
$ws.Range("H7").Value = "Sint-Denijs-Westrem"
$ws.Range("P7").Value = 693300224
$ws.Range("V7").Value = "1907218070029338"

# --- Row 8 (Id 845218) no longer present in the refreshed scrape ---
$ws.Rows(8).Delete()

# Leftover selection from scrolling down to check the table while editing
# the interaction/scraper code.
$ws.Range("D11").Select()
